# Cambios en domiciliaria y PAMI
# Replace the "Instalada" data table (Sede / Servicio / Mes / Capacidad instalada)
# with the new May-2023 figures for Bulevar, San Martin and Cartagena, and grow
# the backing Excel table (Tabla2) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mes is the same for every data row: 2023-05-01 (serial 45047)
$mes = 45047

# New data, in row order (row 2 .. row 41).
# Each entry: Sede, Servicio, Capacidad instalada (or $null when the cell is blank)
$rows = @(
    @("Bulevar",   "MEDICINA GENERAL",               552),
    @("Bulevar",   "PEDIATRÍA",                       300),
    @("Bulevar",   "MEDICINA INTERNA",                346),
    @("Bulevar",   "DERMATOLOGÍA",                    225),
    @("Bulevar",   "PROCEDIMIENTOS DERMATOLÓGICOS",     0),
    @("Bulevar",   "GERIATRÍA",                        30),
    @("Bulevar",   "CARDIOLOGÍA",                     180),
    @("Bulevar",   "ENDOCRINOLOGÍA ",                  60),
    @("Bulevar",   "FISIATRÍA ",                       68),
    @("Bulevar",   "GINECOLOGÍA ",                    228),
    @("Bulevar",   "ECOCARDIOGRAMA ",                 151),
    @("Bulevar",   "ECOGRAFÍA",                       300),
    @("Bulevar",   "NUTRICIÓN",                       360),
    @("Bulevar",   "PSICOLOGÍA",                      360),
    @("Bulevar",   "TRABAJO SOCIAL",                  173),
    @("Bulevar",   "TERAPIA FÍSICA",                 2377),
    @("Bulevar",   "TERAPIA FÍSICA 1RA VEZ",          207),
    @("Bulevar",   "TERAPIA OCUPACIONAL",             368),
    @("Bulevar",   "TERAPIA RESPIRATORIA",            291),
    @("Bulevar",   "TERAPIA DE LENGUAJE",             368),
    @("Bulevar",   "PSIQUIATRÍA",                      40),
    @("Bulevar",   "MEDICINA DEL DEPORTE",             68),
    @("San Martin","MEDICINA GENERAL",                552),
    @("San Martin","MEDICINA INTERNA",                 68),
    @("San Martin","GINECOLOGÍA ",                     48),
    @("San Martin","NUTRICIÓN",                       360),
    @("San Martin","PSICOLOGÍA",                      420),
    @("San Martin","TERAPIA RESPIRATORIA",            $null),
    @("San Martin","TERAPIA FÍSICA",                  344),
    @("San Martin","FISIATRÍA ",                       24),
    @("San Martin","PSIQUIATRÍA",                      14),
    @("Cartagena", "MEDICINA GENERAL",                552),
    @("Cartagena", "TERAPIA FÍSICA",                  572),
    @("Cartagena", "TERAPIA FÍSICA 1RA VEZ",           60),
    @("Cartagena", "MEDICINA INTERNA",                 84),
    @("Cartagena", "NEUROLOGÍA",                       60),
    @("Cartagena", "GERIATRÍA",                        24),
    @("Cartagena", "FISIATRÍA ",                       44),
    @("Cartagena", "REUMATOLOGÍA",                     60),
    @("Cartagena", "PSIQUIATRÍA INFANTIL",             15)
)

# A cell that already carries the short-date display format (style index "s=1"
# in the original file) - used as the source for a format-only copy so every
# new Mes cell gets the same date format without minting a new style.
$dateFormatSource = $ws.Cells.Item(2, 3)

$r = 2
foreach ($row in $rows) {
    $sede = $row[0]
    $servicio = $row[1]
    $capacidad = $row[2]

    $ws.Cells.Item($r, 1).Value = $sede
    $ws.Cells.Item($r, 2).Value = $servicio

    $dateFormatSource.Copy()
    $ws.Cells.Item($r, 3).PasteSpecial(-4122)
    $ws.Cells.Item($r, 3).Value = $mes

    if ($null -eq $capacidad) {
        $ws.Cells.Item($r, 4).ClearContents()
    } else {
        $ws.Cells.Item($r, 4).Value = $capacidad
    }

    $r++
}

# Grow the table (Tabla2) from A1:D37 to A1:D41 so the new rows are included.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D41"))

# Restore the view: select D39 and scroll back to the top (matches the saved
# workbook's sheetView after the edit).
$ws.Range("A1").Select()
$ws.Range("D39").Select()
